$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicadores")

# Update row 7 "Imobilização PL" values
$ws.Range("B7").Value = 1.406442047213015
$ws.Range("C7").Value = 1.429181031520963
$ws.Range("D7").Value = 1.406706186378151
$ws.Range("E7").Value = 1.38777482749692
$ws.Range("F7").Value = 1.373129995073128

# Delete rows 18-23 entirely (GAO, GAF, GAT, Var % Receita, Var % EBIT, Var % Lucro Líquido)
$ws.Range("A18:F23").Delete()
